$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure this sheet/window is the active one so selection changes stick
$ws.Activate()

# Update the header label in A1 ("Hidden Setting" -> "Size of Hidden Layers")
$ws.Range("A1").Value = "Size of Hidden Layers"

# Widen column A to fit the new, longer label (switch from auto best-fit to an explicit custom width)
$ws.Columns.Item(1).ColumnWidth = 12.5

# Update the selected/active cell shown in the sheet view
$ws.Range("F6").Select() | Out-Null
